$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "C.C."
$ws.Range("E2").Value = 1001456
$ws.Range("F2").Value = "CARRERA 50"
$ws.Range("G2").Value = "30-01-2011"
$ws.Range("H2").Value = 14
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 3005656565
$ws.Range("L2").Value = "santiago.@gmail.com"
$ws.Range("M2").Value = "Negocios Internacionales - Virtual"
$ws.Range("O2").Value = "uploaded_files\1001456_CÉDULA.pdf"
$ws.Range("P2").Value = "uploaded_files\1001456_CIVICA.pdf"
$ws.Range("Q2").Value = "uploaded_files\1001456_SERVICIOPUBLICOS.pdf"
$ws.Range("R2").Value = "uploaded_files\1001456_ANEXO1.pdf"
$ws.Range("S2").Value = "uploaded_files\1001456_ANEXO2.xlsx"
